$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 381 (old rows 381.. shift down to 383..)
$ws.Rows.Item(381).Resize(2).Insert()

# New row 381: Pera - Bartlett de verano - Primera
$ws.Range("A381").Value2 = 4
$ws.Range("B381").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C381").Value2 = "Los Lagos"
$ws.Range("D381").Value2 = 44939
$ws.Range("E381").Value2 = 10
$ws.Range("F381").Value2 = "Fruta"
$ws.Range("G381").Value2 = 100104
$ws.Range("H381").Value2 = "Frutos de pepita"
$ws.Range("I381").Value2 = 100104005
$ws.Range("J381").Value2 = "Pera"
$ws.Range("K381").Value2 = "Bartlett de verano"
$ws.Range("L381").Value2 = "Primera"
$ws.Range("M381").Value2 = 400
$ws.Range("N381").Value2 = 22000
$ws.Range("O381").Value2 = 23000
$ws.Range("P381").Value2 = 22500
$ws.Range("Q381").Value2 = "`$/caja 15 kilos empedrada"
$ws.Range("R381").Value2 = "Región de O'Higgins"
$ws.Range("S381").Value2 = 1500
$ws.Range("T381").Value2 = 15

# New row 382: Pera - Bartlett de verano - Segunda
$ws.Range("A382").Value2 = 4
$ws.Range("B382").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C382").Value2 = "Los Lagos"
$ws.Range("D382").Value2 = 44939
$ws.Range("E382").Value2 = 10
$ws.Range("F382").Value2 = "Fruta"
$ws.Range("G382").Value2 = 100104
$ws.Range("H382").Value2 = "Frutos de pepita"
$ws.Range("I382").Value2 = 100104005
$ws.Range("J382").Value2 = "Pera"
$ws.Range("K382").Value2 = "Bartlett de verano"
$ws.Range("L382").Value2 = "Segunda"
$ws.Range("M382").Value2 = 200
$ws.Range("N382").Value2 = 20000
$ws.Range("O382").Value2 = 20000
$ws.Range("P382").Value2 = 20000
$ws.Range("Q382").Value2 = "`$/caja 15 kilos empedrada"
$ws.Range("R382").Value2 = "Región de O'Higgins"
$ws.Range("S382").Value2 = 1333
$ws.Range("T382").Value2 = 15
